$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5899688005447388
$ws.Range("B1").Value = 1.331511616706848
$ws.Range("C1").Value = 5.758915424346924
$ws.Range("D1").Value = 2.594637155532837
$ws.Range("E1").Value = 1.35286271572113
